$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 2.0.0-sd-202312-matchbox-patch -> 2.0.0-sd-202406-matchbox-patch
$meta.Range("B3").Value = "2.0.0-sd-202406-matchbox-patch"

# Date: 2024-03-12T18:28:21+01:00 -> 2024-06-19T17:47:42+02:00
$meta.Range("B8").Value = "2024-06-19T17:47:42+02:00"

# Contact: No display for ContactDetail -> HL7 International - Structured Documents (...)
$meta.Range("B10").Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"

# --- Elements sheet updates ---
$elem = $wb.Worksheets.Item("Elements")

# Binding Value Set (column Z, row 5): v3-SetOperator -> CDASetOperator
$elem.Range("Z5").Value = "http://hl7.org/cda/stds/core/ValueSet/CDASetOperator"

# Update bestFit column width to reflect new (longer) content, matching Excel's
# recalculated best-fit width for the new URL text.
$elem.Columns.Item(26).ColumnWidth = 50.15
